$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(35, 1).Value = "Access-point for resources"
$ws.Cells.Item(35, 2).Value = "Access"
$ws.Cells.Item(35, 3).Value = "/{access-point}?resource={URI}"
$ws.Cells.Item(35, 4).Value = "GET"
$ws.Cells.Item(35, 5).Value = "Empty"
$ws.Cells.Item(35, 6).Value = "Versioning"
$ws.Cells.Item(35, 7).Value = "200 Ok, 404 Not found, 406 Not acceptable"
$ws.Cells.Item(35, 8).Value = "Subject CBD"
$ws.Cells.Item(35, 9).Value = "None"

$ws.Cells.Item(36, 1).Value = "Access-point for resources"
$ws.Cells.Item(36, 2).Value = "Access"
$ws.Cells.Item(36, 3).Value = "/{access-point}?resource={URI}"
$ws.Cells.Item(36, 4).Value = "POST"
$ws.Cells.Item(36, 5).Value = "Form"
$ws.Cells.Item(36, 6).Value = "Versioning"
$ws.Cells.Item(36, 7).Value = "See GET"
$ws.Cells.Item(36, 9).Value = "None"

$ws.Cells.Item(37, 1).Value = "Access-point for resources"
$ws.Cells.Item(37, 2).Value = "Access"
$ws.Cells.Item(37, 3).Value = "/{access-point}?resource={URI}"
$ws.Cells.Item(37, 4).Value = "POST"
$ws.Cells.Item(37, 5).Value = "Content"
$ws.Cells.Item(37, 6).Value = "Context"
$ws.Cells.Item(37, 7).Value = "200 Ok, 400 Bad request, 404 Not found, 409 Conflict"
$ws.Cells.Item(37, 8).Value = "Status info"
$ws.Cells.Item(37, 9).Value = "(version of) content added"

$ws.Cells.Item(38, 1).Value = "Access-point for resources"
$ws.Cells.Item(38, 2).Value = "Access"
$ws.Cells.Item(38, 3).Value = "/{access-point}?resource={URI}"
$ws.Cells.Item(38, 4).Value = "PUT"
$ws.Cells.Item(38, 5).Value = "Content"
$ws.Cells.Item(38, 6).Value = "Context"
$ws.Cells.Item(38, 7).Value = "200 Ok, 400 Bad request, 404 Not found, 409 Conflict"
$ws.Cells.Item(38, 8).Value = "Status info"
$ws.Cells.Item(38, 9).Value = "(version of) content replaced"

$ws.Cells.Item(39, 1).Value = "Access-point for resources"
$ws.Cells.Item(39, 2).Value = "Access"
$ws.Cells.Item(39, 3).Value = "/{access-point}?resource={URI}"
$ws.Cells.Item(39, 4).Value = "DELETE"
$ws.Cells.Item(39, 5).Value = "Empty"
$ws.Cells.Item(39, 6).Value = "Versioning"
$ws.Cells.Item(39, 7).Value = "200 Ok, 400 Bad request, 404 Not found, 409 Conflict"
$ws.Cells.Item(39, 8).Value = "Status info"
$ws.Cells.Item(39, 9).Value = "(version of) content deleted"

$ws.Cells.Item(40, 1).Value = "Access-point for named graphs"
$ws.Cells.Item(40, 2).Value = "Access"
$ws.Cells.Item(40, 3).Value = "/{access-point}?graph={URI}"
$ws.Cells.Item(40, 4).Value = "GET"
$ws.Cells.Item(40, 5).Value = "Empty"
$ws.Cells.Item(40, 6).Value = "Versioning"
$ws.Cells.Item(40, 7).Value = "200 Ok, 404 Not found, 406 Not acceptable"
$ws.Cells.Item(40, 8).Value = "Document content"
$ws.Cells.Item(40, 9).Value = "None"

$ws.Cells.Item(41, 1).Value = "Access-point for named graphs"
$ws.Cells.Item(41, 2).Value = "Access"
$ws.Cells.Item(41, 3).Value = "/{access-point}?graph={URI}"
$ws.Cells.Item(41, 4).Value = "POST"
$ws.Cells.Item(41, 5).Value = "Form"
$ws.Cells.Item(41, 6).Value = "Versioning"
$ws.Cells.Item(41, 7).Value = "See GET"
$ws.Cells.Item(41, 9).Value = "None"

$ws.Cells.Item(42, 1).Value = "Access-point for named graphs"
$ws.Cells.Item(42, 2).Value = "Access"
$ws.Cells.Item(42, 3).Value = "/{access-point}?graph={URI}"
$ws.Cells.Item(42, 4).Value = "POST"
$ws.Cells.Item(42, 5).Value = "Document content"
$ws.Cells.Item(42, 6).Value = "Context"
$ws.Cells.Item(42, 7).Value = "200 Ok, 400 Bad request, 404 Not found, 409 Conflict"
$ws.Cells.Item(42, 8).Value = "Status info"
$ws.Cells.Item(42, 9).Value = "(version of) content added"

$ws.Cells.Item(43, 1).Value = "Access-point for named graphs"
$ws.Cells.Item(43, 2).Value = "Access"
$ws.Cells.Item(43, 3).Value = "/{access-point}?graph={URI}"
$ws.Cells.Item(43, 4).Value = "PUT"
$ws.Cells.Item(43, 5).Value = "Document content"
$ws.Cells.Item(43, 6).Value = "Context"
$ws.Cells.Item(43, 7).Value = "200 Ok, 400 Bad request, 404 Not found, 409 Conflict"
$ws.Cells.Item(43, 8).Value = "Status info"
$ws.Cells.Item(43, 9).Value = "(version of) content replaced"

$ws.Cells.Item(44, 1).Value = "Access-point for named graphs"
$ws.Cells.Item(44, 2).Value = "Access"
$ws.Cells.Item(44, 3).Value = "/{access-point}?graph={URI}"
$ws.Cells.Item(44, 4).Value = "DELETE"
$ws.Cells.Item(44, 5).Value = "Empty"
$ws.Cells.Item(44, 6).Value = "Versioning"
$ws.Cells.Item(44, 7).Value = "200 Ok, 400 Bad request, 404 Not found, 409 Conflict"
$ws.Cells.Item(44, 8).Value = "Status info"
$ws.Cells.Item(44, 9).Value = "(version of) content deleted"


# Set the new column C width so the final OOXML width matches 29
$ws.Columns.Item(3).ColumnWidth = 28.16

# Update selection to mirror the author's final cursor position
$ws.Range("D41").Select()
